# Apply the rubric grading updates described in the commit "updated what I did"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# E10 was marked milestone "II" - update it to milestone "III"
$ws.Range("E10").Value = "III"

# Mark additional rubric items as achieved ("X") for several criteria rows
$ws.Range("F32").Value = "X"
$ws.Range("F34").Value = "X"
$ws.Range("F35").Value = "X"
$ws.Range("F36").Value = "X"
$ws.Range("F40").Value = "X"
$ws.Range("F56").Value = "X"
$ws.Range("F67").Value = "X"

# Update the view state to match where the user was working
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 37
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F61").Select()
